$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared (rich-text) header strings: volume number + report week dates ---
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# --- Cells changing type: number -> existing shared text ("0" / "***.*") ---
# Copy from row 14 template cells (untouched by this edit) to inherit the correct
# style (s=13) + shared string, since these values already exist verbatim.
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("D17"))
$ws.Range("E14").Copy($ws.Range("E17"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("F29"))
$ws.Range("C14").Copy($ws.Range("F30"))

# --- Cells changing type: existing shared text -> number ---
# Copy numeric style template from row 14, then overwrite with the new numeric value.
$ws.Range("I14").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("L14").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("I14").Copy($ws.Range("F31"))
$ws.Range("F31").Value = 2

# --- Plain numeric value updates (style/type unchanged) ---
# Row 16
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 75
$ws.Range("J16").Value = 118
$ws.Range("K16").Value = -36.440677966101
$ws.Range("L16").Value = -41.40625
$ws.Range("M16").Value = -17.582417582417
$ws.Range("N16").Value = -87.13550600343
# Row 17
$ws.Range("C17").Value = 4
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -8.333333333333
$ws.Range("I17").Value = 102
$ws.Range("K17").Value = -12.820512820512
$ws.Range("L17").Value = -32.450331125827
$ws.Range("M17").Value = 59.375
$ws.Range("N17").Value = -52.112676056338
# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -85.714285714285
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -72.222222222222
$ws.Range("I18").Value = 108
$ws.Range("J18").Value = 174
$ws.Range("K18").Value = -37.931034482758
$ws.Range("L18").Value = -45.454545454545
$ws.Range("M18").Value = -16.923076923076
$ws.Range("N18").Value = -81.118881118881
# Row 19
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 93
$ws.Range("G19").Value = 94
$ws.Range("H19").Value = -1.063829787234
$ws.Range("I19").Value = 681
$ws.Range("J19").Value = 753
$ws.Range("K19").Value = -9.561752988047
$ws.Range("L19").Value = -21.634062140391
$ws.Range("M19").Value = -6.584362139917
$ws.Range("N19").Value = -59.415971394517
# Row 20
$ws.Range("I20").Value = 14
$ws.Range("K20").Value = -48.148148148148
$ws.Range("L20").Value = -56.25
$ws.Range("M20").Value = -48.148148148148
$ws.Range("N20").Value = -97.165991902834
# Row 21
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -13.513513513513
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 142
$ws.Range("H21").Value = -14.788732394366
$ws.Range("I21").Value = 986
$ws.Range("J21").Value = 1192
$ws.Range("K21").Value = -17.28187919463
$ws.Range("L21").Value = -28.860028860028
$ws.Range("M21").Value = -5.916030534351
$ws.Range("N21").Value = -72.217526063679
# Row 22
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 32
$ws.Range("K22").Value = 6.666666666666
$ws.Range("L22").Value = 10.344827586206
$ws.Range("M22").Value = -11.111111111111
# Row 24
$ws.Range("C24").Value = 24
$ws.Range("E24").Value = -36.842105263157
$ws.Range("F24").Value = 134
$ws.Range("G24").Value = 176
$ws.Range("H24").Value = -23.863636363636
$ws.Range("I24").Value = 1015
$ws.Range("J24").Value = 1253
$ws.Range("K24").Value = -18.994413407821
$ws.Range("L24").Value = -30.19257221458
$ws.Range("M24").Value = -2.963671128107
# Row 25
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 29
$ws.Range("E25").Value = -48.275862068965
$ws.Range("F25").Value = 99
$ws.Range("G25").Value = 135
$ws.Range("H25").Value = -26.666666666666
$ws.Range("I25").Value = 756
$ws.Range("J25").Value = 996
$ws.Range("K25").Value = -24.096385542168
$ws.Range("L25").Value = -32.19730941704
# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 27
$ws.Range("H26").Value = 17.391304347826
$ws.Range("I26").Value = 252
$ws.Range("J26").Value = 250
$ws.Range("K26").Value = 0.8
$ws.Range("L26").Value = -13.993174061433
$ws.Range("M26").Value = 46.511627906976
# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 53
$ws.Range("K28").Value = 8.163265306122
$ws.Range("L28").Value = 17.777777777777
# Row 31
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 9
$ws.Range("J31").Value = 19
$ws.Range("K31").Value = -52.631578947368
$ws.Range("L31").Value = -18.181818181818

